$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (School category) - replace the old "Question5" question with a new
# question about a friend crying in class.
$ws.Range("C6").Value = "A friend from class is crying, what to do?"
$ws.Range("D6").Value = "She is crying because she's sad, so I will ask her why she's sad"
$ws.Range("E6").Value = "I don't know why she's crying so I'll leave her alone"
$ws.Range("F6").Value = "Laugh because she looks funny"
$ws.Range("H6").Value = "She is crying because she's sad, so I will ask her why she's sad"

# Row 9 (Home category) - update a couple of the answer options for the
# "sister is annoying me" question.
$ws.Range("D9").Value = "Yell at her"
$ws.Range("F9").Value = "Play with her if I can because she wants my attention"

# Row 10 (Home category) - fix the "Right Answer" reference so it points to
# the correct answer text instead of the stray "Answer2" placeholder.
$ws.Range("H10").Value = "Say nothing because it's not polite"

# Row 11 (Home category) - replace the old "Question10" question with a new
# question about a brother hugging.
$ws.Range("C11").Value = "My brother is hugging me but I don't like hugs, what do I do?"
$ws.Range("D11").Value = "Push him away"
$ws.Range("E11").Value = "Yell at him because he's annoying"
$ws.Range("F11").Value = "Tell him that I know he means well but I don't like it"
$ws.Range("H11").Value = "Tell him that I know he means well but I don't like it"
